$d = $word.ActiveDocument

$replacements = @(
    @('2024-01-24 Wednesday', '2024-01-25 Thursday'),
    @('86÷2=43, 0', '70÷3=23, 1'),
    @('11÷4=2, 3', '66÷3=22, 0'),
    @('86÷5=17, 1', '48÷2=24, 0'),
    @('18÷9=2, 0', '96÷7=13, 5'),
    @('85÷8=10, 5', '39÷7=5, 4'),
    @('80÷3=26, 2', '80÷7=11, 3'),
    @('49÷2=24, 1', '86÷7=12, 2'),
    @('43÷5=8, 3', '43÷8=5, 3'),
    @('13÷7=1, 6', '64÷6=10, 4'),
    @('31÷3=10, 1', '32÷3=10, 2'),
    @('33÷8=4, 1', '46÷8=5, 6'),
    @('98÷6=16, 2', '40÷4=10, 0'),
    @('85÷6=14, 1', '64÷4=16, 0'),
    @('83÷6=13, 5', '19÷5=3, 4'),
    @('59÷3=19, 2', '66÷3=22, 0'),
    @('92÷5=18, 2', '38÷4=9, 2'),
    @('29÷7=4, 1', '82÷3=27, 1'),
    @('32÷9=3, 5', '91÷8=11, 3'),
    @('38÷6=6, 2', '52÷5=10, 2'),
    @('60÷4=15, 0', '46÷4=11, 2'),
    @('66÷6=11, 0', '35÷8=4, 3'),
    @('93÷5=18, 3', '45÷3=15, 0'),
    @('16÷2=8, 0', '85÷3=28, 1'),
    @('33÷2=16, 1', '96÷4=24, 0'),
    @('80÷6=13, 2', '57÷6=9, 3'),
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2)
}

$d.Save()